$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.55
$ws.Range("G2").Value = 1.68
$ws.Range("H2").Value = 5.3
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 5.5
$ws.Range("N2").Value = 3.3
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.4
$ws.Range("S2").Value = 2.98
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04
$ws.Range("W2").Value = 2.46
$ws.Range("AN2").Value = 1000
$ws.Range("F3").Value = 6.6
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 1.39
$ws.Range("I3").Value = 1.45
$ws.Range("J3").Value = 1.2
$ws.Range("L3").Value = 1.29
$ws.Range("M3").Value = 1.02
$ws.Range("P3").Value = 2.24
$ws.Range("Q3").Value = 1.59
$ws.Range("R3").Value = 1.54
$ws.Range("S3").Value = 2.28
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.99
$ws.Range("V3").Value = 3.1
$ws.Range("W3").Value = 1.12
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 18
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 21
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 36
$ws.Range("AO3").Value = 1000
$ws.Range("G4").Value = 29
$ws.Range("H4").Value = 1.18
$ws.Range("J4").Value = 1.28
$ws.Range("W4").Value = 1.03
$ws.Range("F5").Value = 1.98
$ws.Range("G5").Value = 2.62
$ws.Range("I5").Value = 4.7
$ws.Range("L5").Value = 1.36
$ws.Range("Q5").Value = 1.93
$ws.Range("S5").Value = 1.94
$ws.Range("V5").Value = 1.27
$ws.Range("W5").Value = 1.62
$ws.Range("F6").Value = 1.77
$ws.Range("G6").Value = 1.9
$ws.Range("H6").Value = 4.6
$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 3
$ws.Range("N6").Value = 1.57
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 1.57
$ws.Range("Q6").Value = 1.42
$ws.Range("R6").Value = 1.18
$ws.Range("S6").Value = 2.4
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.17
$ws.Range("W6").Value = 2.04
$ws.Range("X6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AN6").Value = 1000
